$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "E" (full month forecast) and "G" (confidence interval) columns
# for rows 17-47 (data block), skipping the blank separator rows (23, 30, 40).
$ws.Range("E17").Value = 285314.875
$ws.Range("G17").Value = "266523,125-304106,625"
$ws.Range("E18").Value = 276463.59375
$ws.Range("G18").Value = "252106,59375-300820,59375"
$ws.Range("E19").Value = 273113.5
$ws.Range("G19").Value = "245162,578125-301064,40625"
$ws.Range("E20").Value = 278673.625
$ws.Range("G20").Value = "256724,546875-300622,6875"
$ws.Range("E21").Value = 290103.8125
$ws.Range("G21").Value = "278087,875-302119,75"
$ws.Range("E22").Value = 285718.9375
$ws.Range("G22").Value = "267595,5625-303842,3125"
$ws.Range("E24").Value = 309928.9375
$ws.Range("G24").Value = "288192,5-331665,375"
$ws.Range("E25").Value = 310474.3125
$ws.Range("G25").Value = "288570,4375-332378,1875"
$ws.Range("E26").Value = 310474.3125
$ws.Range("G26").Value = "288570,4375-332378,1875"
$ws.Range("E27").Value = 309351.75
$ws.Range("G27").Value = "286536,125-332167,375"
$ws.Range("E28").Value = 293066.59375
$ws.Range("G28").Value = "281125,375-305007,8125"
$ws.Range("E29").Value = 293043.15625
$ws.Range("G29").Value = "280994,84375-305091,46875"
$ws.Range("E31").Value = 315037.25
$ws.Range("G31").Value = "290692,875-339381,625"
$ws.Range("E32").Value = 314992.71875
$ws.Range("G32").Value = "290443,03125-339542,40625"
$ws.Range("E33").Value = 314992.71875
$ws.Range("G33").Value = "290443,03125-339542,40625"
$ws.Range("E34").Value = 314992.71875
$ws.Range("G34").Value = "290443,03125-339542,40625"
$ws.Range("E35").Value = 314992.71875
$ws.Range("G35").Value = "290443,03125-339542,40625"
$ws.Range("E36").Value = 314992.71875
$ws.Range("G36").Value = "290443,03125-339542,40625"
$ws.Range("E37").Value = 314984.625
$ws.Range("G37").Value = "290429,28125-339539,96875"
$ws.Range("E38").Value = 315104.4375
$ws.Range("G38").Value = "290484,75-339724,125"
$ws.Range("E39").Value = 314775.53125
$ws.Range("G39").Value = "290063,5-339487,5625"
$ws.Range("E41").Value = 241611.703125
$ws.Range("G41").Value = "217422,09375-265801,3125"
$ws.Range("E42").Value = 241791.890625
$ws.Range("G42").Value = "217502,734375-266081,0625"
$ws.Range("E43").Value = 245828.515625
$ws.Range("G43").Value = "219433,21875-272223,8125"
$ws.Range("E44").Value = 242441.703125
$ws.Range("G44").Value = "217789,125-267094,28125"
$ws.Range("E45").Value = 242351.890625
$ws.Range("G45").Value = "217783,9375-266919,84375"
$ws.Range("E46").Value = 242461.84375
$ws.Range("G46").Value = "217761,53125-267162,15625"
$ws.Range("E47").Value = 233130.4375
$ws.Range("G47").Value = "211953,875-254307"

# Update the sheet view: clear the scrolled top-left cell and move the selection
# from the last data row (G47) up to G5.
$null = $ws.Range("G5").Select()
